$d = $word.ActiveDocument

# Locate the target paragraphs by scanning for their distinctive text.
# Two "doSave(order);" paragraphs exist in the document; the one that
# needs editing is immediately followed by the "User DataAccess: ..."
# paragraph, so use that as a disambiguator.
$count = $d.Paragraphs.Count

$target = $null
$loginIdx = $null
$logoutIdx = $null

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.Trim()

    if ($t -eq "doSave(order);") {
        if ($i -lt $count) {
            $next = $d.Paragraphs.Item($i + 1)
            $nextText = $next.Range.Text.Trim()
            if ($nextText.StartsWith("User DataAccess")) {
                $target = $p
            }
        }
    }
    elseif ($t -eq "login();") {
        $loginIdx = $i
    }
    elseif ($t -eq "logout();") {
        $logoutIdx = $i
    }
}

# 1) doSave(order); -> doSave(voceOrdine);
if ($target -ne $null) {
    $target.Range.Find.Execute("order", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "voceOrdine", 2)
}

# 2) Remove the "login();" and "logout();" paragraphs entirely (they sit
#    back-to-back, right after "newCliente(cliente);").
if ($loginIdx -ne $null -and $logoutIdx -ne $null) {
    $pLogin = $d.Paragraphs.Item($loginIdx)
    $pLogout = $d.Paragraphs.Item($logoutIdx)
    $rng = $d.Range($pLogin.Range.Start, $pLogout.Range.End)
    $rng.Delete()
}
